$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on D:E so numeric-looking strings (e.g. "599.61")
# are stored as literal text, matching the original inlineStr cells, not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.886.25'
$ws.Range('E2').Value = '  -0.84%  '

$ws.Range('D3').Value = '3.808.00'
$ws.Range('E3').Value = '  -2.07%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').Value = '599.61'
$ws.Range('E5').Value = '  -0.38%  '

$ws.Range('D6').Value = '168.96'
$ws.Range('E6').Value = '  +0.47%  '

$ws.Range('D7').Value = '3.808.80'
$ws.Range('E7').Value = '  -2.02%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('E9').Value = '  +0.35%  '

$ws.Range('E10').Value = '  -0.90%  '

$ws.Range('D11').Value = '6.50'
$ws.Range('E11').Value = '  +0.74%  '

$ws.Range('D12').Value = '0.463'
$ws.Range('E12').Value = '  +0.75%  '

$ws.Range('D13').Value = '0.0000276'
$ws.Range('E13').Value = '  +9.83%  '

$ws.Range('D14').Value = '36.94'
$ws.Range('E14').Value = '  -0.54%  '

$ws.Range('D15').Value = '4.449.30'
$ws.Range('E15').Value = '  -1.93%  '

$ws.Range('D16').Value = '3.796.72'
$ws.Range('E16').Value = '  -2.06%  '

$ws.Range('D17').Value = '67.971.79'
$ws.Range('E17').Value = '  -0.55%  '

$ws.Range('D18').Value = '18.47'
$ws.Range('E18').Value = '  +1.47%  '

$ws.Range('D19').Value = '7.42'
$ws.Range('E19').Value = '  -0.12%  '

$ws.Range('E20').Value = '  +0.30%  '

$ws.Range('D21').Value = '10.87'
$ws.Range('E21').Value = '  +0.00%  '

$ws.Range('D22').Value = '469.66'
$ws.Range('E22').Value = '  -0.79%  '

$ws.Range('D23').Value = '0.737'
$ws.Range('E23').Value = '  -0.01%  '

$ws.Range('D24').Value = '0.0000151'
$ws.Range('E24').Value = '  -8.98%  '

$ws.Range('D25').Value = '83.36'
$ws.Range('E25').Value = '  -0.58%  '

$ws.Range('D26').Value = '2.33'
$ws.Range('E26').Value = '  +2.98%  '

$ws.Range('D27').Value = '12.26'
$ws.Range('E27').Value = '  +0.22%  '

$ws.Range('D28').Value = '10.29'
$ws.Range('E28').Value = '  +2.47%  '

$ws.Range('E29').Value = '  -0.09%  '

$ws.Range('D30').Value = '2.93'
$ws.Range('E30').Value = '  -1.31%  '

$ws.Range('D31').Value = '3.953.36'
$ws.Range('E31').Value = '  -2.09%  '

$ws.Range('D32').Value = '7.75'
$ws.Range('E32').Value = '  -2.26%  '

$ws.Range('D33').Value = '2.29'
$ws.Range('E33').Value = '  -1.31%  '

$ws.Range('D34').Value = '30.80'
$ws.Range('E34').Value = '  -2.33%  '

$ws.Range('D35').Value = '9.35'
$ws.Range('E35').Value = '  -0.82%  '

$ws.Range('D36').Value = '3.774.29'
$ws.Range('E36').Value = '  -2.30%  '

$ws.Range('E37').Value = '  +1.95%  '

$ws.Range('D38').Value = '3.77'
$ws.Range('E38').Value = '  +6.16%  '

$ws.Range('D39').Value = '5.99'
$ws.Range('E39').Value = '  +0.98%  '

$ws.Range('E40').Value = '  -0.92%  '

$ws.Range('E41').Value = '  -2.07%  '

$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.12%  '

$ws.Range('D43').Value = '0.320'
$ws.Range('E43').Value = '  +1.63%  '

$ws.Range('D45').Value = '8.80'
$ws.Range('E45').Value = '  +2.10%  '

$ws.Range('D46').Value = '1.98'
$ws.Range('E46').Value = '  -1.11%  '

$ws.Range('D47').Value = '410.54'
$ws.Range('E47').Value = '  -5.34%  '

$ws.Range('B48').Value = 'FLOKI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D48').Value = '0.000287'
$ws.Range('E48').Value = '  -1.38%  '

$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '46.43'
$ws.Range('E49').Value = '  -1.89%  '

$ws.Range('D50').Value = '142.28'
$ws.Range('E50').Value = '  -1.05%  '

# Restore default cell style (removes the temporary Text number format again)
$ws.Range("D2:E51").Style = "Normal"

Write-Output "done"